$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$sub5 = [char]0x2085
$sub8 = [char]0x2088

Set-TextValue $ws.Range("B2") "Bitcoin"
Set-TextValue $ws.Range("C2") "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue $ws.Range("D2") "26.097.24"
Set-TextValue $ws.Range("E2") "  -0.82%  "

Set-TextValue $ws.Range("B3") "Ethereum"
Set-TextValue $ws.Range("C3") "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue $ws.Range("D3") "1.652.45"
Set-TextValue $ws.Range("E3") "  -0.89%  "

Set-TextValue $ws.Range("B4") "TetherUSD"
Set-TextValue $ws.Range("C4") "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue $ws.Range("D4") "1.004"
Set-TextValue $ws.Range("E4") "  -0.44%  "

Set-TextValue $ws.Range("B5") "BNB"
Set-TextValue $ws.Range("C5") "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue $ws.Range("D5") "218.63"
Set-TextValue $ws.Range("E5") "  -0.77%  "

Set-TextValue $ws.Range("B6") "XRP"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D6") "0.5254"
Set-TextValue $ws.Range("E6") "  -1.02%  "

Set-TextValue $ws.Range("B7") "USDC"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.004"
Set-TextValue $ws.Range("E7") "  -0.43%  "

Set-TextValue $ws.Range("B8") "Cardano"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D8") "0.2670"
Set-TextValue $ws.Range("E8") "  +0.90%  "

Set-TextValue $ws.Range("B9") "Dogecoin"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.06367"
Set-TextValue $ws.Range("E9") "  +0.13%  "

Set-TextValue $ws.Range("B10") "Solana"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D10") "20.56"
Set-TextValue $ws.Range("E10") "  -1.78%  "

Set-TextValue $ws.Range("B11") "TRON"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D11") "0.07685"
Set-TextValue $ws.Range("E11") "  -1.97%  "

Set-TextValue $ws.Range("B12") "Polkadot"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D12") "4.594"
Set-TextValue $ws.Range("E12") "  +1.50%  "

Set-TextValue $ws.Range("B13") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "1.880.15"
Set-TextValue $ws.Range("E13") "  -0.80%  "

Set-TextValue $ws.Range("B14") "Polygon"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D14") "0.5609"
Set-TextValue $ws.Range("E14") "  +0.05%  "

Set-TextValue $ws.Range("B15") "WrappedEther"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "1.487.36"
Set-TextValue $ws.Range("E15") "  -10.86%  "

Set-TextValue $ws.Range("B16") "ShibaInu"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$dVal = "0.0{0}8228" -f $sub5
Set-TextValue $ws.Range("D16") $dVal
Set-TextValue $ws.Range("E16") "  +1.22%  "

Set-TextValue $ws.Range("B17") "Litecoin"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "65.46"
Set-TextValue $ws.Range("E17") "  -0.41%  "

Set-TextValue $ws.Range("B18") "WrappedBTC"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "26.100.62"
Set-TextValue $ws.Range("E18") "  -0.77%  "

Set-TextValue $ws.Range("B19") "Dai"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D19") "1.004"
Set-TextValue $ws.Range("E19") "  -0.47%  "

Set-TextValue $ws.Range("B20") "Uniswap"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "4.689"
Set-TextValue $ws.Range("E20") "  -0.46%  "

Set-TextValue $ws.Range("B21") "Avalanche"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D21") "10.36"
Set-TextValue $ws.Range("E21") "  +0.81%  "

Set-TextValue $ws.Range("B22") "BitcoinCash"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D22") "191.10"
Set-TextValue $ws.Range("E22") "  -3.95%  "

Set-TextValue $ws.Range("B23") "Chainlink"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D23") "5.971"
Set-TextValue $ws.Range("E23") "  -1.36%  "

Set-TextValue $ws.Range("B24") "BinanceUSD"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D24") "1.005"
Set-TextValue $ws.Range("E24") "  -0.43%  "

Set-TextValue $ws.Range("B25") "Monero"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D25") "146.03"
Set-TextValue $ws.Range("E25") "  -0.44%  "

Set-TextValue $ws.Range("B26") "Stellar"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D26") "0.1201"
Set-TextValue $ws.Range("E26") "  -1.08%  "

Set-TextValue $ws.Range("B27") "Cosmos"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "7.250"
Set-TextValue $ws.Range("E27") "  +0.17%  "

Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "15.94"
Set-TextValue $ws.Range("E28") "  -1.60%  "

Set-TextValue $ws.Range("B29") "Toncoin"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D29") "1.497"
Set-TextValue $ws.Range("E29") "  -0.95%  "

Set-TextValue $ws.Range("B30") "Hedera"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.05638"
Set-TextValue $ws.Range("E30") "  -4.48%  "

Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.271"
Set-TextValue $ws.Range("E31") "  -1.01%  "

Set-TextValue $ws.Range("B32") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D32") "3.498"
Set-TextValue $ws.Range("E32") "  -0.82%  "

Set-TextValue $ws.Range("B33") "Filecoin"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D33") "3.379"
Set-TextValue $ws.Range("E33") "  +1.91%  "

Set-TextValue $ws.Range("B34") "LidoDAOToken"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D34") "1.579"
Set-TextValue $ws.Range("E34") "  -1.39%  "

Set-TextValue $ws.Range("B35") "MXToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D35") "2.798"
Set-TextValue $ws.Range("E35") "  -1.03%  "

Set-TextValue $ws.Range("B36") "ARBITRUM"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "0.9452"
Set-TextValue $ws.Range("E36") "  -1.63%  "

Set-TextValue $ws.Range("B37") "HuobiToken"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D37") "2.407"
Set-TextValue $ws.Range("E37") "  -0.98%  "

Set-TextValue $ws.Range("B38") "ImmutableX"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.5776"
Set-TextValue $ws.Range("E38") "  -0.46%  "

Set-TextValue $ws.Range("B39") "VeChain"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.01593"
Set-TextValue $ws.Range("E39") "  -1.42%  "

Set-TextValue $ws.Range("B40") "FraxShare"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D40") "5.969"
Set-TextValue $ws.Range("E40") "  +0.14%  "

Set-TextValue $ws.Range("B41") "PaxDollar"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D41") "1.004"
Set-TextValue $ws.Range("E41") "  -0.48%  "

Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "0.8388"
Set-TextValue $ws.Range("E42") "  -2.15%  "

Set-TextValue $ws.Range("B43") "Maker"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D43") "1.022.52"
Set-TextValue $ws.Range("E43") "  -4.89%  "

Set-TextValue $ws.Range("B44") "Quant"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D44") "101.51"
Set-TextValue $ws.Range("E44") "  -1.24%  "

Set-TextValue $ws.Range("B45") "RocketPoolETH"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D45") "1.791.65"
Set-TextValue $ws.Range("E45") "  -0.78%  "

Set-TextValue $ws.Range("B46") "Aave"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "58.58"
Set-TextValue $ws.Range("E46") "  +0.28%  "

Set-TextValue $ws.Range("B47") "Cronos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.05343"
Set-TextValue $ws.Range("E47") "  +3.80%  "

Set-TextValue $ws.Range("B48") "Frax"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D48") "1.003"
Set-TextValue $ws.Range("E48") "  -0.95%  "

Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.047"
Set-TextValue $ws.Range("E49") "  -0.39%  "

Set-TextValue $ws.Range("B50") "Mantle"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D50") "0.4343"
Set-TextValue $ws.Range("E50") "  -1.58%  "

Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$dVal = "0.0{0}103" -f $sub8
Set-TextValue $ws.Range("D51") $dVal
Set-TextValue $ws.Range("E51") "  -0.37%  "
